{"js": "// Update stack-trace line numbers in the document body text to reflect the\n// M2DocEvaluator.java / AbstractTemplatesTestSuite.java line shifts caused by\n// moving from M2Doc 2.0.0 to 2.0.1.\nconst replacements = [\n  [\"M2DocEvaluator.java:1049)\", \"M2DocEvaluator.java:1061)\"],\n  [\"M2DocEvaluator.java:1084)\", \"M2DocEvaluator.java:1096)\"],\n  [\"M2DocEvaluator.java:1300)\", \"M2DocEvaluator.java:1305)\"],\n  [\"M2DocEvaluator.java:278)\", \"M2DocEvaluator.java:283)\"],\n  [\"M2DocEvaluator.java:267)\", \"M2DocEvaluator.java:272)\"],\n  [\"AbstractTemplatesTestSuite.java:475)\", \"AbstractTemplatesTestSuite.java:479)\"],\n  [\"AbstractTemplatesTestSuite.java:384)\", \"AbstractTemplatesTestSuite.java:388)\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // Each old string may appear more than once (e.g. the \"1084\" -> \"1096\"\n  // line repeats 3 times), so keep searching/replacing until no more hits.\n  // eslint-disable-next-line no-constant-condition\n  while (true) {\n    const results = body.search(oldText, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      break;\n    }\n\n    for (const range of results.items) {\n      range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Update stack-trace line numbers in the document body text to reflect the\n# M2DocEvaluator.java / AbstractTemplatesTestSuite.java line shifts caused by\n# moving from M2Doc 2.0.0 to 2.0.1.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"M2DocEvaluator.java:1049)\", \"M2DocEvaluator.java:1061)\"),\n    @(\"M2DocEvaluator.java:1084)\", \"M2DocEvaluator.java:1096)\"),\n    @(\"M2DocEvaluator.java:1300)\", \"M2DocEvaluator.java:1305)\"),\n    @(\"M2DocEvaluator.java:278)\", \"M2DocEvaluator.java:283)\"),\n    @(\"M2DocEvaluator.java:267)\", \"M2DocEvaluator.java:272)\"),\n    @(\"AbstractTemplatesTestSuite.java:475)\", \"AbstractTemplatesTestSuite.java:479)\"),\n    @(\"AbstractTemplatesTestSuite.java:384)\", \"AbstractTemplatesTestSuite.java:388)\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
